$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 159
$ws.Range("J2").Value = 54.50153350830078

# Row 3
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 159
$ws.Range("J3").Value = 186.4371299743652

# Row 4
$ws.Range("C4").Value = "sac"
$ws.Range("D4").Value = "suc"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 160
$ws.Range("J4").Value = 191.0860538482666

# Row 5
$ws.Range("C5").Value = "suc"
$ws.Range("D5").Value = "sac"
$ws.Range("F5").Value = 160
$ws.Range("J5").Value = 193.4034824371338

# Row 6
$ws.Range("F6").Value = 160
$ws.Range("J6").Value = 193.6280727386475

# Row 7
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 160
$ws.Range("J7").Value = 196.1033344268799

# Row 8
$ws.Range("F8").Value = 160
$ws.Range("J8").Value = 195.2712535858154

# Row 9
$ws.Range("F9").Value = 160
$ws.Range("J9").Value = 193.2954788208008
